# Auto-generated edit script: updates betting-odds cell values per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 2.7
$ws.Range("I3").Value = 2.75
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("Z3").Value = 26
$ws.Range("AO3").Value = 15
# Row 4
$ws.Range("X4").Value = 9
$ws.Range("Z4").Value = 15
$ws.Range("AC4").Value = 13
$ws.Range("AE4").Value = 13
$ws.Range("AG4").Value = 151
$ws.Range("AP4").Value = 17
$ws.Range("AU4").Value = 7.5
# Row 5
$ws.Range("O5").Value = 1.22
$ws.Range("P5").Value = 4.33
$ws.Range("Q5").Value = 1.73
$ws.Range("R5").Value = 2.1
# Row 7
$ws.Range("G7").Value = 1.38
$ws.Range("H7").Value = 4.75
$ws.Range("I7").Value = 8
$ws.Range("J7").Value = 1.83
$ws.Range("N7").Value = 15
$ws.Range("W7").Value = 9
$ws.Range("Y7").Value = 9
$ws.Range("Z7").Value = 9.5
$ws.Range("AQ7").Value = 17
# Row 9
$ws.Range("H9").Value = 3.25
$ws.Range("I9").Value = 1.85
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 2.6
$ws.Range("M9").Value = 1.1
$ws.Range("N9").Value = 7
$ws.Range("Y9").Value = 15
$ws.Range("AC9").Value = 7
$ws.Range("AE9").Value = 19
$ws.Range("AF9").Value = 67
$ws.Range("AJ9").Value = 9
$ws.Range("AM9").Value = 34
$ws.Range("AU9").Value = 9
$ws.Range("AV9").Value = 67
$ws.Range("BC9").Value = 201
# Row 11
$ws.Range("U11").Value = 1.44
$ws.Range("V11").Value = 2.63
# Row 12
$ws.Range("G12").Value = 2.45
$ws.Range("I12").Value = 2.75
$ws.Range("J12").Value = 3
$ws.Range("L12").Value = 3.25
$ws.Range("Q12").Value = 1.7
$ws.Range("R12").Value = 2.1
$ws.Range("U12").Value = 1.57
$ws.Range("W12").Value = 11
$ws.Range("Y12").Value = 10
$ws.Range("Z12").Value = 23
$ws.Range("AA12").Value = 19
$ws.Range("AD12").Value = 6.5
$ws.Range("AI12").Value = 15
$ws.Range("AN12").Value = 4.75
$ws.Range("AO12").Value = 13
$ws.Range("AS12").Value = 126
$ws.Range("BA12").Value = 41
# Row 13
$ws.Range("G13").Value = 2.1
$ws.Range("I13").Value = 3.2
$ws.Range("J13").Value = 2.88
$ws.Range("M13").Value = 1.05
$ws.Range("N13").Value = 11
$ws.Range("O13").Value = 1.3
$ws.Range("P13").Value = 3.4
$ws.Range("Q13").Value = 2.03
$ws.Range("R13").Value = 1.83
$ws.Range("S13").Value = 1.4
$ws.Range("T13").Value = 2.75
$ws.Range("U13").Value = 1.8
$ws.Range("V13").Value = 1.95
$ws.Range("W13").Value = 7.5
$ws.Range("X13").Value = 10
$ws.Range("AG13").Value = 251
$ws.Range("AK13").Value = 34
$ws.Range("AL13").Value = 26
$ws.Range("AM13").Value = 34
$ws.Range("AN13").Value = 4.33
$ws.Range("AT13").Value = 2.75
$ws.Range("AX13").Value = 5
$ws.Range("AZ13").Value = 26
$ws.Range("BA13").Value = 51
# Row 14
$ws.Range("G14").Value = 2.5
$ws.Range("I14").Value = 2.63
$ws.Range("J14").Value = 3.25
$ws.Range("L14").Value = 3.5
$ws.Range("N14").Value = 9
$ws.Range("U14").Value = 1.91
$ws.Range("V14").Value = 1.91
$ws.Range("X14").Value = 12
$ws.Range("Y14").Value = 10
$ws.Range("Z14").Value = 26
$ws.Range("AD14").Value = 6.5
$ws.Range("AH14").Value = 8
$ws.Range("AL14").Value = 23
$ws.Range("AN14").Value = 4.5
$ws.Range("AO14").Value = 15
$ws.Range("AR14").Value = 81
$ws.Range("AZ14").Value = 26
# Row 15
$ws.Range("G15").Value = 1.25
$ws.Range("H15").Value = 5.5
$ws.Range("I15").Value = 9
$ws.Range("J15").Value = 1.67
$ws.Range("L15").Value = 8
$ws.Range("O15").Value = 1.13
$ws.Range("P15").Value = 6
$ws.Range("Q15").Value = 1.44
$ws.Range("R15").Value = 2.7
$ws.Range("W15").Value = 9.5
$ws.Range("X15").Value = 7
$ws.Range("Y15").Value = 9.5
$ws.Range("AB15").Value = 26
$ws.Range("AH15").Value = 23
$ws.Range("AJ15").Value = 23
$ws.Range("AO15").Value = 5.5
$ws.Range("AQ15").Value = 13
# Row 17
$ws.Range("G17").Value = 1.53
$ws.Range("H17").Value = 4.1
$ws.Range("I17").Value = 5.75
$ws.Range("J17").Value = 2.1
$ws.Range("K17").Value = 2.25
$ws.Range("L17").Value = 6
$ws.Range("S17").Value = 1.36
$ws.Range("T17").Value = 3
$ws.Range("U17").Value = 2
$ws.Range("V17").Value = 1.73
$ws.Range("X17").Value = 7
$ws.Range("Z17").Value = 11
$ws.Range("AC17").Value = 11
$ws.Range("AD17").Value = 8
$ws.Range("AF17").Value = 67
$ws.Range("AH17").Value = 15
$ws.Range("AI17").Value = 29
$ws.Range("AJ17").Value = 19
$ws.Range("AK17").Value = 67
$ws.Range("AM17").Value = 51
$ws.Range("AN17").Value = 3.5
$ws.Range("AO17").Value = 8
$ws.Range("AP17").Value = 19
$ws.Range("AQ17").Value = 23
$ws.Range("AT17").Value = 3
$ws.Range("AU17").Value = 9
$ws.Range("AX17").Value = 7.5
$ws.Range("AY17").Value = 34
$ws.Range("AZ17").Value = 41
$ws.Range("BA17").Value = 126
$ws.Range("BB17").Value = 151
$ws.Range("BC17").Value = 301
# Row 20
$ws.Range("G20").Value = 4.05
$ws.Range("H20").Value = 3.65
$ws.Range("I20").Value = 1.72
$ws.Range("J20").Value = 4.5
$ws.Range("K20").Value = 2.18
$ws.Range("L20").Value = 2.32
$ws.Range("N20").Value = 7.5
$ws.Range("P20").Value = 3.25
$ws.Range("Q20").Value = 1.9
$ws.Range("R20").Value = 1.85
$ws.Range("S20").Value = 1.39
$ws.Range("T20").Value = 2.77
$ws.Range("U20").Value = 1.85
$ws.Range("V20").Value = 1.87
$ws.Range("W20").Value = 11.5
$ws.Range("X20").Value = 22
$ws.Range("Y20").Value = 14
$ws.Range("Z20").Value = 65
$ws.Range("AA20").Value = 40
$ws.Range("AC20").Value = 7.5
$ws.Range("AD20").Value = 7.1
$ws.Range("AE20").Value = 16.5
$ws.Range("AF20").Value = 80
$ws.Range("AI20").Value = 8
$ws.Range("AJ20").Value = 8.25
$ws.Range("AK20").Value = 13.5
$ws.Range("AL20").Value = 14
$ws.Range("AN20").Value = 5.9
$ws.Range("AO20").Value = 23
$ws.Range("AS20").Value = 450
$ws.Range("AT20").Value = 2.77
$ws.Range("AU20").Value = 7.7
$ws.Range("AX20").Value = 3.6
$ws.Range("AY20").Value = 8.5
$ws.Range("AZ20").Value = 18.5
$ws.Range("BA20").Value = 30
$ws.Range("BB20").Value = 65
# Row 22
$ws.Range("G22").Value = 3.4
$ws.Range("K22").Value = 2.05
$ws.Range("W22").Value = 9
$ws.Range("AE22").Value = 15
$ws.Range("AG22").Value = 351
$ws.Range("AH22").Value = 7
$ws.Range("AK22").Value = 21
$ws.Range("AN22").Value = 5
$ws.Range("AV22").Value = 51
$ws.Range("AZ22").Value = 23
# Row 24
$ws.Range("G24").Value = 1.53
$ws.Range("H24").Value = 3.8
$ws.Range("I24").Value = 6.5
$ws.Range("M24").Value = 1.06
$ws.Range("N24").Value = 10
$ws.Range("Q24").Value = 1.93
$ws.Range("R24").Value = 1.93
$ws.Range("AC24").Value = 10
$ws.Range("AD24").Value = 7.5
$ws.Range("AE24").Value = 19
$ws.Range("AG24").Value = 351
$ws.Range("AI24").Value = 29
$ws.Range("AL24").Value = 41
$ws.Range("AX24").Value = 7.5
$ws.Range("AY24").Value = 34
$ws.Range("BA24").Value = 126
$ws.Range("BB24").Value = 151
# Row 25
$ws.Range("G25").Value = 1.91
$ws.Range("H25").Value = 3.4
$ws.Range("I25").Value = 3.75
$ws.Range("J25").Value = 2.6
$ws.Range("O25").Value = 1.25
$ws.Range("P25").Value = 3.75
$ws.Range("Q25").Value = 1.85
$ws.Range("R25").Value = 1.95
$ws.Range("X25").Value = 9.5
$ws.Range("AI25").Value = 19
$ws.Range("AM25").Value = 34
$ws.Range("AX25").Value = 5.5
